$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Company")

# Row 12 = "D&B Rating" (under Corporate Viability). The scorer didn't get
# to rate it, so the numeric score is cleared and a note explains why.
$ws.Range("B12").ClearContents()

$notes = $ws.Range("C12")
$notes.Value = "No time to check"
# Match the formatting already used by the sheet's other populated "Notes"
# cells (bordered, bold Arial) rather than the blank-cell default.
$notes.VerticalAlignment = -4107
